$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 '28.191.20'
Set-TextValue 2 5 '  -0.31%  '

# Row 3
Set-TextValue 3 4 '1.910.55'
Set-TextValue 3 5 '  +2.03%  '

# Row 4
Set-TextValue 4 5 '  -0.11%  '

# Row 5
Set-TextValue 5 4 '314.76'

# Row 6
Set-TextValue 6 5 '  -0.11%  '

# Row 7
Set-TextValue 7 5 '  +0.66%  '

# Row 8
Set-TextValue 8 4 '0.3929'
Set-TextValue 8 5 '  -0.01%  '

# Row 9
Set-TextValue 9 4 '0.09323'
Set-TextValue 9 5 '  -3.11%  '

# Row 10
Set-TextValue 10 4 '1.140'
Set-TextValue 10 5 '  -0.61%  '

# Row 11
Set-TextValue 11 4 '41.89'
Set-TextValue 11 5 '  +2.41%  '

# Row 12
Set-TextValue 12 4 '6.405'
Set-TextValue 12 5 '  -1.22%  '

# Row 13
Set-TextValue 13 5 '  -0.56%  '

# Row 14
Set-TextValue 14 4 '1.900.83'
Set-TextValue 14 5 '  +1.20%  '

# Row 15
Set-TextValue 15 4 '7.325'
Set-TextValue 15 5 '  -1.32%  '

# Row 16
Set-TextValue 16 5 '  -0.14%  '

# Row 17
Set-TextValue 17 4 '0.00001123'
Set-TextValue 17 5 '  -0.63%  '

# Row 18
Set-TextValue 18 4 '92.43'
Set-TextValue 18 5 '  -0.46%  '

# Row 19
Set-TextValue 19 4 '0.06617'

# Row 20
Set-TextValue 20 4 '17.98'
Set-TextValue 20 5 '  +1.84%  '

# Row 21
Set-TextValue 21 4 '1.000'
Set-TextValue 21 5 '  -0.08%  '

# Row 22
Set-TextValue 22 4 '6.234'
Set-TextValue 22 5 '  +0.65%  '

# Row 23
Set-TextValue 23 4 '28.243.32'
Set-TextValue 23 5 '  -0.34%  '

# Row 24
Set-TextValue 24 5 '  +1.78%  '

# Row 25
Set-TextValue 25 4 '2.323'
Set-TextValue 25 5 '  +1.45%  '

# Row 26
Set-TextValue 26 4 '2.594'
Set-TextValue 26 5 '  +0.99%  '

# Row 27
Set-TextValue 27 4 '2.124.63'
Set-TextValue 27 5 '  +1.48%  '

# Row 28
Set-TextValue 28 4 '21.11'
Set-TextValue 28 5 '  -0.70%  '

# Row 29
Set-TextValue 29 4 '158.03'
Set-TextValue 29 5 '  -0.47%  '

# Row 30
Set-TextValue 30 4 '127.22'
Set-TextValue 30 5 '  -0.25%  '

# Row 31
Set-TextValue 31 5 '  +3.48%  '

# Row 32
Set-TextValue 32 5 '  +0.91%  '

# Row 33
Set-TextValue 33 4 '5.651'
Set-TextValue 33 5 '  +0.38%  '

# Row 34
Set-TextValue 34 5 '  -0.23%  '

# Row 35
Set-TextValue 35 4 '9.708'
Set-TextValue 35 5 '  +2.14%  '

# Row 36
Set-TextValue 36 4 '0.06684'
Set-TextValue 36 5 '  -0.94%  '

# Row 37
Set-TextValue 37 4 '0.02427'
Set-TextValue 37 5 '  +1.09%  '

# Row 38
Set-TextValue 38 4 '1.243'
Set-TextValue 38 5 '  -0.15%  '

# Row 39
Set-TextValue 39 4 '0.2203'
Set-TextValue 39 5 '  +0.58%  '

# Row 40
Set-TextValue 40 4 '1.287'
Set-TextValue 40 5 '  +8.62%  '

# Row 41
Set-TextValue 41 4 '0.6516'
Set-TextValue 41 5 '  +2.58%  '

# Row 42
Set-TextValue 42 4 '11.53'
Set-TextValue 42 5 '  +0.22%  '

# Row 43
Set-TextValue 43 4 '5.009'
Set-TextValue 43 5 '  +0.12%  '

# Row 44
Set-TextValue 44 4 '1.000'
Set-TextValue 44 5 '  -0.08%  '

# Row 45
Set-TextValue 45 2 'Decentraland'
Set-TextValue 45 3 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 45 4 '0.6098'
Set-TextValue 45 5 '  +1.80%  '

# Row 46
Set-TextValue 46 2 'EnergySwap'
Set-TextValue 46 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 46 4 '13.36'
Set-TextValue 46 5 '  -1.40%  '

# Row 47
Set-TextValue 47 5 '  +1.70%  '

# Row 48
Set-TextValue 48 4 '1.289'
Set-TextValue 48 5 '  +1.37%  '

# Row 49
Set-TextValue 49 4 '2.020'
Set-TextValue 49 5 '  +0.81%  '

# Row 50
Set-TextValue 50 4 '123.63'
Set-TextValue 50 5 '  -0.50%  '

# Row 51
Set-TextValue 51 4 '1.189'
Set-TextValue 51 5 '  -0.58%  '
